$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 12-16 (A:C) are reshuffled: "Dru Smith" is removed and
# "Brandon Boston Jr." is added, with Rudy Gobert / Dejounte Murray /
# Jalen Suggs moving to new rows.

$ws.Range("A12").Value = "Rudy Gobert"
$ws.Range("B12").Value = "C"
$ws.Range("C12").Value = "Minnesota Timberwolves"

$ws.Range("A13").Value = "Malcolm Brogdon"
$ws.Range("B13").Value = "PG,SG"
$ws.Range("C13").Value = "Washington Wizards"

$ws.Range("A14").Value = "Brandon Boston Jr."
$ws.Range("B14").Value = "SG,SF,PF"
$ws.Range("C14").Value = "New Orleans Pelicans"

$ws.Range("A15").Value = "Dejounte Murray"
$ws.Range("B15").Value = "PG,SG"
$ws.Range("C15").Value = "New Orleans Pelicans"

$ws.Range("A16").Value = "Jalen Suggs"
$ws.Range("B16").Value = "PG,SG"
$ws.Range("C16").Value = "Orlando Magic"
